$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Delete column A - shifts B:F left to A:E
$ws.Columns("A").Delete()

# Rename the shared-string header text from MODEL_CONDITION to MODELCONDITION
# After the column deletion, the header that used to read "MODEL_CONDITION" is now in column D (row 1).
$ws.Cells.Item(1, 4).Value = "MODELCONDITION"
